# Apply the changes described in the diff to sheet1 of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header row to the english/short machine-friendly names.
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Title-case the "de"/"del"/"los" connector words inside a handful of
#    place names (state / municipality columns).
$ws.Range("B2").Value = "San José De Gracia"
$ws.Range("A23").Value = "Ciudad De México"
$ws.Range("A36").Value = "Estado De México"
$ws.Range("B36").Value = "Almoloya De Juárez"
$ws.Range("B40").Value = "Ecatepec De Morelos"
$ws.Range("B53").Value = "Acapulco De Juárez"
$ws.Range("B56").Value = "Atenango Del Río"
$ws.Range("B58").Value = "Chilapa De Álvarez"
$ws.Range("B59").Value = "Chilpancingo De Los Bravo"
$ws.Range("B63").Value = "Tlapa De Comonfort"
$ws.Range("B68").Value = "Lagos De Moreno"
$ws.Range("B78").Value = "Tlaltizapán De Zapata"
$ws.Range("B84").Value = "Chalcatongo De Hidalgo"
$ws.Range("B85").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B87").Value = "Ixtlán De Juárez"
$ws.Range("B90").Value = "Mariscala De Juárez"
$ws.Range("B91").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B103").Value = "Tataltepec De Valdés"
$ws.Range("B113").Value = "Izúcar De Matamoros"
$ws.Range("B120").Value = "Tepanco De López"

# 3. Drop the trailing footnote rows (152-156): sample size, source,
#    author and date notes. This also shrinks the sheet dimension down
#    to A1:D150.
$ws.Rows("152:156").Delete()
